# Commit: "added 50x reliability data"
#
# Adds a new worksheet "50x 100uL" right after "All Data", holding 50
# trial measurements for a single 100uL reading (column A = trial number
# 1..50 via a running "+1" formula, column B = measured weight).

$wb = $excel.ActiveWorkbook

$allData = $wb.Worksheets.Item("All Data")

# Insert the new worksheet immediately after "All Data" -- this becomes
# the active sheet, which also clears tabSelected on "All Data" and sets
# the workbook's active tab to this new sheet (index 1).
$ws = $wb.Worksheets.Add($null, $allData)
$ws.Name = "50x 100uL"

# Headers
$ws.Cells.Item(1, 1).Value = "Trial"
$ws.Cells.Item(1, 2).Value = "100uL"

# 50 trial rows: B column is the measured value, A column is the trial
# index (1 literal, then a running =previous+1 formula down to 50).
$values = @(
    97.8,
    98,
    98,
    98.5,
    97.5,
    97.5,
    97.2,
    96.9,
    96.2,
    97.6,
    97.2,
    96.6,
    96.6,
    99.7,
    95.8,
    95.7,
    100.3,
    96,
    96,
    97.8,
    95,
    99.3,
    95.1,
    96.7,
    95.9,
    96.3,
    97.4,
    98.9,
    97.2,
    96.3,
    96.1,
    95.4,
    95.5,
    95.1,
    95.5,
    96.1,
    95.5,
    96.2,
    95.1,
    95.7,
    95.9,
    95.4,
    95.2,
    96.2,
    99.1,
    97.2,
    97.1,
    98.6,
    98.4,
    94.1
)

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = $values[0]

for ($i = 1; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $prevRow = $row - 1
    $ws.Cells.Item($row, 1).Formula = "=A$prevRow+1"
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Match the author's final view state: scrolled down with D48 selected.
$ws.Range("D48").Select()
